$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the row above (row 9) into the new row 10, A column
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A10").Value = "2021年"

$ws.Range("B10").Value = 40
$ws.Range("C10").Value = 34.4
$ws.Range("D10").Value = 104.7
$ws.Range("E10").Value = 63.8
$ws.Range("F10").Value = 75.7
$ws.Range("G10").Value = 36.7
$ws.Range("H10").Value = 35.4
$ws.Range("I10").Value = 36.2
$ws.Range("J10").Value = 34.4
$ws.Range("K10").Value = 72.40000000000001
$ws.Range("L10").Value = 52.2
$ws.Range("M10").Value = 48
$ws.Range("N10").Value = 49.6
$ws.Range("O10").Value = 46.4
$ws.Range("P10").Value = 69.8
$ws.Range("Q10").Value = 49.6
$ws.Range("R10").Value = 30.5
